# Weekly fruit/vegetable price update.
# Rows 442-444 (old data) are copied down to new rows 445-447 (append),
# then rows 441-444 are overwritten in place with the new weekly figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$dateFmt = "YYYY-MM-DD HH:MM:SS"

# --- 1. Append the pre-existing rows 442, 443, 444 as new rows 445, 446, 447 ---

# New row 445 = old row 442
$ws.Range("A445").Value = 4
$ws.Range("B445").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C445").Value = "Los Lagos"
$ws.Range("D445").Value = 44552
$ws.Range("D445").NumberFormat = $dateFmt
$ws.Range("E445").Value = 10
$ws.Range("F445").Value = 100112033
$ws.Range("G445").Value = "Lechuga"
$ws.Range("H445").Value = "Escarola"
$ws.Range("I445").Value = "Primera"
$ws.Range("J445").Value = 150
$ws.Range("K445").Value = 10000
$ws.Range("L445").Value = 10000
$ws.Range("M445").Value = 10000
$ws.Range("N445").Value = "`$/caja 15 unidades"
$ws.Range("O445").Value = "Región Metropolitana"
$ws.Range("P445").Value = 667
$ws.Range("Q445").Value = 15
$ws.Range("R445").Value = "Hortaliza"

# New row 446 = old row 443
$ws.Range("A446").Value = 4
$ws.Range("B446").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C446").Value = "Los Lagos"
$ws.Range("D446").Value = 44544
$ws.Range("D446").NumberFormat = $dateFmt
$ws.Range("E446").Value = 10
$ws.Range("F446").Value = 100112033
$ws.Range("G446").Value = "Lechuga"
$ws.Range("H446").Value = "Escarola"
$ws.Range("I446").Value = "Primera"
$ws.Range("J446").Value = 500
$ws.Range("K446").Value = 7500
$ws.Range("L446").Value = 8000
$ws.Range("M446").Value = 7750
$ws.Range("N446").Value = "`$/caja 15 unidades"
$ws.Range("O446").Value = "Región Metropolitana"
$ws.Range("P446").Value = 517
$ws.Range("Q446").Value = 15
$ws.Range("R446").Value = "Hortaliza"

# New row 447 = old row 444
$ws.Range("A447").Value = 4
$ws.Range("B447").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C447").Value = "Los Lagos"
$ws.Range("D447").Value = 44160
$ws.Range("D447").NumberFormat = $dateFmt
$ws.Range("E447").Value = 10
$ws.Range("F447").Value = 100112033
$ws.Range("G447").Value = "Lechuga"
$ws.Range("H447").Value = "Escarola"
$ws.Range("I447").Value = "Primera"
$ws.Range("J447").Value = 40
$ws.Range("K447").Value = 6000
$ws.Range("L447").Value = 6000
$ws.Range("M447").Value = 6000
$ws.Range("N447").Value = "`$/caja 15 unidades"
$ws.Range("O447").Value = "Región de Coquimbo"
$ws.Range("P447").Value = 400
$ws.Range("Q447").Value = 15
$ws.Range("R447").Value = "Hortaliza"

# --- 2. Overwrite rows 441-444 in place with the new weekly values ---

# Row 441
$ws.Range("D441").Value = 44595
$ws.Range("H441").Value = "Conconina(o)"
$ws.Range("J441").Value = 120
$ws.Range("K441").Value = 11000
$ws.Range("L441").Value = 11000
$ws.Range("M441").Value = 11000
$ws.Range("N441").Value = "`$/caja 10 unidades"
$ws.Range("O441").Value = "Región Metropolitana"
$ws.Range("P441").Value = 1100
$ws.Range("Q441").Value = 10

# Row 442
$ws.Range("D442").Value = 44595
$ws.Range("J442").Value = 100
$ws.Range("K442").Value = 13000
$ws.Range("L442").Value = 13000
$ws.Range("M442").Value = 13000
$ws.Range("O442").Value = "Región de Coquimbo"
$ws.Range("P442").Value = 867

# Row 443
$ws.Range("D443").Value = 44595
$ws.Range("I443").Value = "Segunda"
$ws.Range("J443").Value = 100
$ws.Range("K443").Value = 11000
$ws.Range("L443").Value = 11000
$ws.Range("M443").Value = 11000
$ws.Range("N443").Value = "`$/caja 18 unidades"
$ws.Range("O443").Value = "Región de Coquimbo"
$ws.Range("P443").Value = 611
$ws.Range("Q443").Value = 18

# Row 444
$ws.Range("D444").Value = 44335
$ws.Range("J444").Value = 100
$ws.Range("K444").Value = 12000
$ws.Range("L444").Value = 12000
$ws.Range("M444").Value = 12000
$ws.Range("P444").Value = 800
